$wb = $excel.ActiveWorkbook

$headers = @("variable","var.ratio","bias","cor1","cor2","RMSE")

function Set-SheetData {
    param($ws, $rows)
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
    }
    $r = 2
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        for ($c = 1; $c -lt $row.Length; $c++) {
            $ws.Cells.Item($r, $c + 1).Value = $row[$c]
        }
        $r++
    }
}

# lm sheet (index 1)
$wsLm = $wb.Worksheets.Item(1)
$LmRows = @(
    @("no", 0.356998431242402, 1.00875263869728, 0.552086240708944, 0.426337282112785, 0.0844028528077384),
    @("no2", 0.39596828970968, 1.01509176833617, 0.568361087528583, 0.517560291044462, 0.114226974056548),
    @("o3", 0.623553116802343, 1.00774907911811, 0.740450608175213, 0.483941005009877, 0.122890541132252),
    @("pm10", 0.249214539534343, 1.01142299048137, 0.431418413758017, 0.402512323377383, 0.0700250975727749),
    @("pm2.5", 0.300304688464977, 1.0049538872016, 0.501817263513791, 0.494460855914543, 0.0629765123633122)
)
Set-SheetData $wsLm $LmRows

# glm sheet (index 2)
$wsGlm = $wb.Worksheets.Item(2)
$GlmRows = @(
    @("no", 0.355682851764593, 0.999757505513779, 0.566697826252613, 0.429571041318111, 0.0831003921894197),
    @("no2", 0.397401740429509, 1.00067369285668, 0.574277788695762, 0.51907014242039, 0.113641500380643),
    @("o3", 0.628840182120175, 1.00133507253272, 0.728161871076884, 0.482633404897965, 0.125864862063917),
    @("pm10", 0.249144669442368, 1.00013433680816, 0.434031265766093, 0.403684969393035, 0.0699160014350932),
    @("pm2.5", 0.298298918013315, 0.999724705744873, 0.506789188959009, 0.498322764429342, 0.0627125444755543)
)
Set-SheetData $wsGlm $GlmRows

# randomForest sheet (index 3)
$wsRandomForest = $wb.Worksheets.Item(3)
$RandomForestRows = @(
    @("no", 0.459824848727384, 1.02436926363346, 0.631902385956943, 0.348674713822861, 0.078140142844817),
    @("no2", 0.469121674351911, 1.01517019750327, 0.633623020824343, 0.461580994381967, 0.108656374507281),
    @("o3", 0.630481516356916, 1.00563755265042, 0.764373174498036, 0.409791018860439, 0.114464054233356),
    @("pm10", 0.371148283548642, 1.01727770806536, 0.464223784349427, 0.322631624997902, 0.0675322428447305),
    @("pm2.5", 0.411292790334665, 1.01031335079449, 0.610790164026597, 0.424258698899956, 0.0585189603330702)
)
Set-SheetData $wsRandomForest $RandomForestRows

# knn sheet (new, added after randomForest)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsKnn = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsKnn.Name = "knn"
$knnRows = @(
    @("no", 0.993061188896763, 0.999069111035801, 0.446463898863761, 0.166378922887274, 0.105559552304538),
    @("no2", 1.06714774953197, 1.00339570979271, 0.430246451282144, 0.276068453452492, 0.147158163469892),
    @("o3", 1.02401062559034, 1.00783036493845, 0.661847778321802, 0.296779377150589, 0.150204845491451),
    @("pm10", 1.57816846432379, 0.998758505285383, 0.274939944402421, 0.135988749341796, 0.0934634978501901),
    @("pm2.5", 1.22979895225807, 1.00754215929452, 0.395971697514185, 0.207192293035423, 0.0821613982076782)
)
Set-SheetData $wsKnn $knnRows

# Restore the originally active sheet/tab (lm) so sheet selection state
# matches the pre-edit workbook instead of leaving the newly added sheet active.
$wsLm.Activate()

